$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '69.821.36'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +5.13%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.647.67'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +18.50%  '

$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.09%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '592.36'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +3.08%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '185.84'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +8.93%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '3.650.67'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +18.72%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.999'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -0.14%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.535'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +4.90%  '

$ws.Range('E10').Value = '  +7.37%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '6.53'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +4.13%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.496'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +5.53%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '39.27'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +9.92%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.0000253'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +5.93%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '4.244.45'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +18.23%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '3.638.86'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +18.27%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '69.847.62'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +5.33%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.124'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +1.88%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '7.55'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +8.48%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '17.14'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +3.37%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '509.03'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +4.91%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '9.35'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +21.64%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.745'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +8.63%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '87.82'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +6.68%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '13.51'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +6.76%  '

$ws.Range('E26').Value = '  +8.01%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '10.86'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +7.18%  '

$ws.Range('E28').Value = '  -0.04%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.54'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +12.60%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '8.17'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +2.91%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '32.65'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +17.65%  '

$ws.Range('E32').Value = '  +5.48%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.0000108'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +18.85%  '

$ws.Range('E34').Value = '  +5.52%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.999'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -0.03%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '6.14'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +10.53%  '

$ws.Range('E37').Value = '  +8.97%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.335'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +11.17%  '

$ws.Range('E39').Value = '  +7.70%  '

$ws.Range('B40').Value = 'OKB'
$ws.Range('C40').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '50.86'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +3.71%  '

$ws.Range('B41').Value = 'Arweave'
$ws.Range('C41').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '46.53'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -3.67%  '

$ws.Range('E42').Value = '  +4.41%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '3.157.97'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +14.07%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '8.84'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +7.44%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '2.78'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +10.05%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '404.24'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +10.72%  '

$ws.Range('E47').Value = '  +6.84%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '27.80'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +14.57%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '136.90'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +1.69%  '

$ws.Range('E50').Value = '  +0.07%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '2.45'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +13.94%  '
